# JDB TPA3 du 06/02/2024
# Journal de bord update: correct the teacher code on row 10 (M -> CM),
# and fill in the new TPA3 session row (row 12).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 (2024-02-05, FSIL #9): the "Enseignant" column had an inconsistent
# value "M" which should have been "CM".
$ws.Range("C10").Value = "CM"

# Row 12: new session entry (2024-02-06), MPAL, TP, group A3,
# with description and comments about the QuestionsScore kata session (TPA3).
$ws.Range("A12").Value = 45328
$ws.Range("B12").Value = "MPAL"
$ws.Range("C12").Value = "TP"
$ws.Range("F12").Value = "x"
$ws.Range("G12").Value = "QuestionsScore Fix #1 en kata"
$ws.Range("I12").Value = "On a avancé jusqu'à écrire des cas nominaux et cas limites pour le constructeur de QCE (énoncé vide ou null, indice < 1), et du score pour indice < 1. Factorisation avec BeforeEach"

# Leave the cursor on H12, matching the author's last edited position.
$ws.Range("H12").Select()
